# Auto-generated script applying market-data refresh values to Leve profit sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 2144.7273
$ws.Range("I62").Value = 1209.7778
$ws.Range("J62").Value = 2792
$ws.Range("K62").Value = 1209.7778
$ws.Range("L62").Value = 2792
$ws.Range("M62").Value = -585.7778000000001
$ws.Range("N62").Value = -4040
$ws.Range("H65").Value = 2144.7273
$ws.Range("I65").Value = 1209.7778
$ws.Range("J65").Value = 2792
$ws.Range("K65").Value = 6048.889
$ws.Range("L65").Value = 13960
$ws.Range("M65").Value = -2928.889
$ws.Range("N65").Value = -20200
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()
$ws.Range("H98").Value = 896.1111
$ws.Range("I98").Value = 772.35297
$ws.Range("J98").Value = 3000
$ws.Range("K98").Value = 772.35297
$ws.Range("L98").Value = 3000
$ws.Range("M98").Value = 725.64703
$ws.Range("N98").Value = -5996
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()
$ws.Range("H112").Value = 2791.923
$ws.Range("J112").Value = 2932.9167
$ws.Range("L112").Value = 8798.750100000001
$ws.Range("N112").Value = -11014.7501
$ws.Range("H122").Value = 896.1111
$ws.Range("I122").Value = 772.35297
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 2317.05891
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = 132.9410899999998
$ws.Range("N122").Value = -13900
$ws.Range("H123").Value = 30000
$ws.Range("J123").Value = 30000
$ws.Range("L123").Value = 30000
$ws.Range("N123").Value = -39800
$ws.Range("H138").Value = 2307.8774
$ws.Range("J138").Value = 3311.7
$ws.Range("L138").Value = 9935.099999999999
$ws.Range("N138").Value = -20215.1

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H12").Value = 1000
$ws.Range("J12").Value = 1000
$ws.Range("L12").Value = 1000
$ws.Range("N12").Value = -1346
$ws.Range("H31").Value = 17586.834
$ws.Range("I31").Value = 2104.2
$ws.Range("K31").Value = 2104.2
$ws.Range("M31").Value = -1810.2
$ws.Range("H61").Value = 5378142
$ws.Range("I61").Value = 7752978
$ws.Range("J61").Value = 3512.7896
$ws.Range("K61").Value = 7752978
$ws.Range("L61").Value = 3512.7896
$ws.Range("M61").Value = -7752766
$ws.Range("N61").Value = -3936.7896
$ws.Range("H74").Value = 1040.8182
$ws.Range("I74").Value = 632.3570999999999
$ws.Range("J74").Value = 1341.7894
$ws.Range("K74").Value = 632.3570999999999
$ws.Range("L74").Value = 1341.7894
$ws.Range("M74").Value = 241.6429000000001
$ws.Range("N74").Value = -3089.7894
$ws.Range("H77").Value = 1040.8182
$ws.Range("I77").Value = 632.3570999999999
$ws.Range("J77").Value = 1341.7894
$ws.Range("K77").Value = 3161.7855
$ws.Range("L77").Value = 6708.946999999999
$ws.Range("M77").Value = 1206.2145
$ws.Range("N77").Value = -15444.947
$ws.Range("H93").Value = 70224
$ws.Range("J93").Value = 70224
$ws.Range("L93").Value = 70224
$ws.Range("N93").Value = -75216
$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()
$ws.Range("H127").Value = 42780
$ws.Range("J127").Value = 42780
$ws.Range("L127").Value = 42780
$ws.Range("N127").Value = -52700
$ws.Range("H132").Value = 4813.469
$ws.Range("I132").Value = 4863.9116
$ws.Range("K132").Value = 14591.7348
$ws.Range("M132").Value = -12061.7348
$ws.Range("H134").Value = 69926.10000000001
$ws.Range("J134").Value = 69926.10000000001
$ws.Range("L134").Value = 69926.10000000001
$ws.Range("N134").Value = -80066.10000000001
$ws.Range("H136").Value = 5378142
$ws.Range("I136").Value = 7752978
$ws.Range("J136").Value = 3512.7896
$ws.Range("K136").Value = 23258934
$ws.Range("L136").Value = 10538.3688
$ws.Range("M136").Value = -23256384
$ws.Range("N136").Value = -15638.3688
$ws.Range("H141").Value = 50425.668
$ws.Range("J141").Value = 50425.668
$ws.Range("L141").Value = 50425.668
$ws.Range("N141").Value = -60785.668

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H13").Value = 70356
$ws.Range("J13").Value = 70356
$ws.Range("L13").Value = 70356
$ws.Range("N13").Value = -70692
$ws.Range("H51").Value = 94000
$ws.Range("J51").Value = 94000
$ws.Range("L51").Value = 94000
$ws.Range("N51").Value = -94982
$ws.Range("H52").Value = 94000
$ws.Range("J52").Value = 94000
$ws.Range("L52").Value = 94000
$ws.Range("N52").Value = -94526
$ws.Range("H55").Value = 47277.5
$ws.Range("J55").Value = 47277.5
$ws.Range("L55").Value = 47277.5
$ws.Range("N55").Value = -47823.5
$ws.Range("H59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").ClearContents()
$ws.Range("H96").Value = 32476
$ws.Range("I96").Value = 32476
$ws.Range("K96").Value = 32476
$ws.Range("M96").Value = -29730
$ws.Range("H121").Value = 94000
$ws.Range("J121").Value = 94000
$ws.Range("L121").Value = 94000
$ws.Range("N121").Value = -97494
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
$ws.Range("H134").Value = 3119.7666
$ws.Range("I134").Value = 3119.5881
$ws.Range("J134").Value = 3120
$ws.Range("K134").Value = 9358.764299999999
$ws.Range("L134").Value = 9360
$ws.Range("M134").Value = -6823.764299999999
$ws.Range("N134").Value = -14430

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 48332.668
$ws.Range("J20").Value = 48332.668
$ws.Range("L20").Value = 48332.668
$ws.Range("N20").Value = -48804.668
$ws.Range("H30").Value = 48332.668
$ws.Range("J30").Value = 48332.668
$ws.Range("L30").Value = 48332.668
$ws.Range("N30").Value = -48514.668
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()
$ws.Range("H128").Value = 48332.668
$ws.Range("J128").Value = 48332.668
$ws.Range("L128").Value = 48332.668
$ws.Range("N128").Value = -58292.668

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H35").Value = 4772.2856
$ws.Range("J35").Value = 5517.6665
$ws.Range("L35").Value = 16552.9995
$ws.Range("N35").Value = -17128.9995
$ws.Range("H68").Value = 1269.8334
$ws.Range("J68").Value = 1471.9722
$ws.Range("L68").Value = 4415.9166
$ws.Range("N68").Value = -6037.9166
$ws.Range("H71").Value = 1269.8334
$ws.Range("J71").Value = 1471.9722
$ws.Range("L71").Value = 13247.7498
$ws.Range("N71").Value = -21359.7498
$ws.Range("H137").Value = 9208
$ws.Range("J137").Value = 2933.3333
$ws.Range("L137").Value = 8799.999899999999
$ws.Range("N137").Value = -18999.9999
$ws.Range("H140").Value = 1531.2162
$ws.Range("I140").Value = 1058.0714
$ws.Range("J140").Value = 3003.2222
$ws.Range("K140").Value = 3174.2142
$ws.Range("L140").Value = 9009.6666
$ws.Range("M140").Value = 2005.7858
$ws.Range("N140").Value = -19369.6666

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 65002
$ws.Range("I12").Value = 60000
$ws.Range("K12").Value = 60000
$ws.Range("M12").Value = -59860
$ws.Range("H62").Value = 30000
$ws.Range("J62").Value = 30000
$ws.Range("L62").Value = 30000
$ws.Range("N62").Value = -31372
$ws.Range("H65").Value = 30000
$ws.Range("J65").Value = 30000
$ws.Range("L65").Value = 90000
$ws.Range("N65").Value = -96864
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H19").Value = 10000000
$ws.Range("I19").Value = 10000000
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 10000000
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = -9999830
$ws.Range("N19").ClearContents()
$ws.Range("H70").Value = 87500
$ws.Range("J70").Value = 87500
$ws.Range("L70").Value = 87500
$ws.Range("N70").Value = -88040
$ws.Range("H73").Value = 87500
$ws.Range("J73").Value = 87500
$ws.Range("L73").Value = 87500
$ws.Range("N73").Value = -89372
$ws.Range("H102").Value = 33900
$ws.Range("J102").Value = 33900
$ws.Range("L102").Value = 33900
$ws.Range("N102").Value = -40390
$ws.Range("H123").Value = 30000
$ws.Range("J123").Value = 30000
$ws.Range("L123").Value = 30000
$ws.Range("N123").Value = -39800
$ws.Range("H140").Value = 48429
$ws.Range("J140").Value = 48429
$ws.Range("L140").Value = 48429
$ws.Range("N140").Value = -58789

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H42").Value = 70035.664
$ws.Range("J42").Value = 70035.664
$ws.Range("L42").Value = 70035.664
$ws.Range("N42").Value = -70791.664
$ws.Range("H111").Value = 94321.5
$ws.Range("J111").Value = 94321.5
$ws.Range("L111").Value = 94321.5
$ws.Range("N111").Value = -102501.5
$ws.Range("H129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("N129").ClearContents()
